# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# list with freshly scraped values (GitHub Actions refresh).
#
# The sheet stores these as plain text cells (not numbers), since the
# "Price" strings use "." as both a thousands- and decimal-separator
# (e.g. "60.939.56") and would be mangled/misread if stored numerically.
# For price values that Excel COM would otherwise auto-coerce to a
# number (e.g. "1.00" -> 1), a leading apostrophe forces text entry,
# exactly like typing into a cell in the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.939.56'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '3.381.86'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("D5").Value = '''571.02'
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").Value = '''141.78'
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '''7.60'
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("E10").Value = '  -1.45%  '
$ws.Range("D11").Value = '''0.388'
$ws.Range("E11").Value = '  -1.50%  '
$ws.Range("D12").Value = '3.962.09'
$ws.Range("E12").Value = '  -0.19%  '
$ws.Range("E13").Value = '  +1.96%  '
$ws.Range("D14").Value = '''27.85'
$ws.Range("E14").Value = '  -1.44%  '
$ws.Range("E15").Value = '  -0.04%  '
$ws.Range("D16").Value = '3.395.84'
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").Value = '61.049.38'
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").Value = '''6.07'
$ws.Range("E18").Value = '  -3.59%  '
$ws.Range("E19").Value = '  -3.90%  '
$ws.Range("D20").Value = '''8.96'
$ws.Range("E20").Value = '  -1.64%  '
$ws.Range("D21").Value = '''382.93'
$ws.Range("D22").Value = '''75.09'
$ws.Range("E22").Value = '  +2.72%  '
$ws.Range("D23").Value = '''0.550'
$ws.Range("E23").Value = '  -2.54%  '
$ws.Range("E24").Value = '  +0.40%  '
$ws.Range("E25").Value = '  -2.60%  '
$ws.Range("D26").Value = '3.522.27'
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").Value = '''7.25'
$ws.Range("E29").Value = '  -1.32%  '
$ws.Range("D30").Value = '''7.97'
$ws.Range("E30").Value = '  -2.43%  '
$ws.Range("D31").Value = '''2.15'
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("D33").Value = '''1.38'
$ws.Range("E33").Value = '  -4.27%  '
$ws.Range("D34").Value = '''23.21'
$ws.Range("E34").Value = '  -2.48%  '
$ws.Range("D35").Value = '''6.93'
$ws.Range("E35").Value = '  -0.64%  '
$ws.Range("D36").Value = '''166.82'
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("D37").Value = '3.414.12'
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("D38").Value = '''4.97'
$ws.Range("E38").Value = '  -2.10%  '
$ws.Range("E39").Value = '  -4.40%  '
$ws.Range("D40").Value = '''0.0768'
$ws.Range("D41").Value = '''26.86'
$ws.Range("E41").Value = '  -1.05%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("E43").Value = '  -1.16%  '
$ws.Range("D44").Value = '''4.37'
$ws.Range("E44").Value = '  -2.60%  '
$ws.Range("D45").Value = '''1.66'
$ws.Range("E45").Value = '  -2.33%  '
$ws.Range("E46").Value = '  -0.28%  '
$ws.Range("D47").Value = '2.449.92'
$ws.Range("E47").Value = '  -3.82%  '
$ws.Range("D48").Value = '''22.95'
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("E49").Value = '  -2.61%  '
$ws.Range("D50").Value = '''2.15'
$ws.Range("E50").Value = '  +8.34%  '
$ws.Range("E51").Value = '  +1.33%  '
